# MAJ 20/02/2018
# - refresh the "last updated" date footer (14/02/2018 -> 19/02/2018) on the
#   slide master and every slide layout
# - append a new, empty "Titre et contenu" slide at the end of the deck

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached datetimeFigureOut footer field everywhere it appears:
#    the slide master and all 11 custom layouts. We reach every layout by
#    briefly adding a throwaway slide on that layout (CustomLayout.Shapes is
#    only reliable when walked via a Slide that uses it), patch the date
#    placeholder text on the shared layout part, then drop the scratch slide.
# ---------------------------------------------------------------------------

function Update-DateShape($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        $txt = ""
        try { $txt = $shp.TextFrame.TextRange.Text } catch {}
        if ($txt -eq "14/02/2018") {
            $shp.TextFrame.TextRange.Text = "19/02/2018"
        }
    }
}

$layoutIds = @(1, 2, 33, 4, 5, 6, 7, 8, 9, 10, 27)
foreach ($layoutId in $layoutIds) {
    $scratch = $p.Slides.Add($p.Slides.Count + 1, $layoutId)
    Update-DateShape $scratch.CustomLayout.Shapes
    $scratch.Delete()
}

Update-DateShape $p.Slides.Item(1).Master.Shapes

# ---------------------------------------------------------------------------
# 2) Append the new slide 18 (Title and Content layout, left empty).
# ---------------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

Write-Output "slides=$($p.Slides.Count) newLayout=$($newSlide.CustomLayout.Name)"
